# Apply updates described in the commit diff (gh-pages data refresh):
#
# Sheet "展览" (exhibitions):
#   F4  22   -> 24    (11.16合肥动漫同人only)
#   F5  4126 -> 4221  (第九届环形宇宙动漫游戏嘉年华)
#   F8  259  -> 261   (心动恋章·冬日序国乙&代号鹄同人only)
#   F9  26   -> 30    (蓬莱次元青年文化节)
#
# Sheet "演出" (performances):
#   F2  123 -> 124    (《四月是你的谎言》音乐会)
#
# Sheet "全部类型" (all types, aggregated):
#   F3  123  -> 124   (《四月是你的谎言》音乐会)
#   F8  22   -> 24    (11.16合肥动漫同人only)
#   F9  4126 -> 4221  (第九届环形宇宙动漫游戏嘉年华)
#   F13 259  -> 261   (心动恋章·冬日序国乙&代号鹄同人only)
#   F14 26   -> 30    (蓬莱次元青年文化节)

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 24
$wsExhibit.Range("F5").Value = 4221
$wsExhibit.Range("F8").Value = 261
$wsExhibit.Range("F9").Value = 30

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 124

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 124
$wsAll.Range("F8").Value = 24
$wsAll.Range("F9").Value = 4221
$wsAll.Range("F13").Value = 261
$wsAll.Range("F14").Value = 30
